$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Insert a new worksheet "Svetla" between "Optika" and "Kamery"
# ---------------------------------------------------------------
$kontrolery = $wb.Worksheets.Item("Kontrolery")
$kamery = $wb.Worksheets.Item("Kamery")

$ws = $wb.Worksheets.Add()
$ws.Name = "Světla"

# Column A width (matches source sheet's first data column)
$ws.Columns.Item(1).ColumnWidth = 24

# Merge the header cells first, then copy the header-cell formatting
# from the "Kontrolery" sheet (merging a freshly-formatted range makes
# Excel rewrite the shared border edge, which would give a distinct
# style from the source; merging first avoids that).
$ws.Range("A1:B1").Merge()
$ws.Range("C1:D1").Merge()

$kontrolery.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)

$kontrolery.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("A1").Value = "všichni výrobci"

$values = @("CL-50","BL-50W-4S","LA-70B","LA-70W","LA-120W","LL-130W","DLU-140W-HI","DL-120W-HI","SP-27IR-850","FL-MD90MC")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Move the new sheet so it sits right before "Kamery" (i.e. right
# after "Optika"), matching the target tab order.
$ws.Move($kamery)

# Re-resolve the sheet by name (the reference can go stale across a
# Move) and make it the active / selected sheet+cell.
$ws = $wb.Worksheets.Item("Světla")
$ws.Activate()
$ws.Range("A1").Select()
